$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133, shifting existing rows 133-198 down to 134-199
$ws.Rows.Item(133).Insert()

# Populate the new row 133 with the new record's data
$ws.Range("A133").Value2 = 11
$ws.Range("B133").Value2 = "Vega Monumental Concepción"
$ws.Range("C133").Value2 = "Bíobío"
$ws.Range("D133").Value2 = 44636
$ws.Range("E133").Value2 = 8
$ws.Range("F133").Value2 = 100114013
$ws.Range("G133").Value2 = "Zanahoria"
$ws.Range("H133").Value2 = "Sin especificar"
$ws.Range("I133").Value2 = "Primera"
$ws.Range("J133").Value2 = 220
$ws.Range("K133").Value2 = 6000
$ws.Range("L133").Value2 = 6500
$ws.Range("M133").Value2 = 6227
$ws.Range("N133").Value2 = "$/saco 20 kilos"
$ws.Range("O133").Value2 = "Chillán"
$ws.Range("P133").Value2 = 311
$ws.Range("Q133").Value2 = 20
$ws.Range("R133").Value2 = "Hortaliza"
